$d = $word.ActiveDocument

# --- Locate the paragraph whose text is exactly "Anexos" -------------------
$paras = $d.Paragraphs
$count = $paras.Count
$targetIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $paras.Item($i)
    $t = $p.Range.Text.Trim()
    if ($t -eq "Anexos") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find the 'Anexos' paragraph"
}

$anexos = $d.Paragraphs.Item($targetIndex)

# --- Create a placeholder paragraph right after "Anexos" -------------------
# (InsertParagraphAfter splits the paragraph mark, producing a new empty
#  paragraph positioned exactly between "Anexos" and whatever followed it.)
$anexos.Range.InsertParagraphAfter()

# Re-fetch the freshly created placeholder paragraph (the previous handle is
# stale after the structural edit above).
$placeholder = $d.Paragraphs.Item($targetIndex + 1)

$w = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$xml = @"
<w:p $w>
  <w:pPr>
    <w:pStyle w:val="Normal"/>
    <w:widowControl w:val="false"/>
    <w:spacing w:before="120" w:after="0"/>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:eastAsia="Arial" w:cs="Calibri"/>
      <w:b/>
      <w:b/>
      <w:i/>
      <w:i/>
      <w:iCs/>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:eastAsia="Arial" w:cs="Calibri" w:ascii="Calibri" w:hAnsi="Calibri"/>
      <w:b/>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
  </w:r>
</w:p>
<w:p $w>
  <w:pPr>
    <w:pStyle w:val="Normal"/>
    <w:widowControl w:val="false"/>
    <w:spacing w:before="120" w:after="0"/>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:eastAsia="Arial" w:cs="Calibri"/>
      <w:b/>
      <w:b/>
      <w:i/>
      <w:i/>
      <w:iCs/>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:eastAsia="Arial" w:cs="Calibri" w:ascii="Calibri" w:hAnsi="Calibri"/>
      <w:b w:val="false"/>
      <w:bCs w:val="false"/>
      <w:i/>
      <w:iCs/>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
    <w:t>Download do APK:</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:eastAsia="Arial" w:cs="Calibri" w:ascii="Calibri" w:hAnsi="Calibri"/>
      <w:b/>
      <w:i/>
      <w:iCs/>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
    <w:br/>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:eastAsia="Arial" w:cs="Calibri" w:ascii="Calibri" w:hAnsi="Calibri"/>
      <w:b w:val="false"/>
      <w:bCs w:val="false"/>
      <w:i w:val="false"/>
      <w:iCs w:val="false"/>
      <w:color w:val="3465A4"/>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
    <w:t>https://bresodev.github.io/aprendIA/APK/aprendIA.apk</w:t>
  </w:r>
</w:p>
"@

$placeholder.Range.InsertXML($xml)
